$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 336, pushing existing rows 336-354 down to 337-355.
$ws.Rows(336).Insert()

# Populate the new row 336 with the new weekly data point.
$ws.Range("A336").Value = 10
$ws.Range("B336").Value = "Vega Modelo de Temuco"
$ws.Range("C336").Value = "La Araucanía"
$ws.Range("D336").Value = 44706
$ws.Range("E336").Value = 9
$ws.Range("F336").Value = 100112040
$ws.Range("G336").Value = "Cilantro"
$ws.Range("H336").Value = "Sin especificar"
$ws.Range("I336").Value = "Primera"
$ws.Range("J336").Value = 80
$ws.Range("K336").Value = 3300
$ws.Range("L336").Value = 3300
$ws.Range("M336").Value = 3300
$ws.Range("N336").Value = "$/docena de atados (2 kilos)"
$ws.Range("O336").Value = "Región Metropolitana"
$ws.Range("P336").Value = 1650
$ws.Range("Q336").Value = 2
$ws.Range("R336").Value = "Hortaliza"

# Match the date formatting used by the rest of column D.
$ws.Range("D336").NumberFormat = $ws.Range("D337").NumberFormat
